{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the async (context) => { ... } function content.\n\n// Helper: replace the first occurrence of `before` text with `after` text,\n// using an exact, case-sensitive, non-wildcard search against the body.\nasync function replaceOnce(context, before, after) {\n  const results = context.document.body.search(before, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + before);\n  }\n\n  results.items[0].insertText(after, \"Replace\");\n  await context.sync();\n}\n\n// 1. Ativa\u00e7\u00e3o date bump.\nawait replaceOnce(\n  context,\n  \"Ativa\u00e7\u00e3o: 01/01/2022\",\n  \"Ativa\u00e7\u00e3o: 01/01/2024\"\n);\n\n// 2. Curso (semestre ideal): drop the EQD entry.\nawait replaceOnce(\n  context,\n  \"Curso (semestre ideal): EQD (3), EQN (4)\",\n  \"Curso (semestre ideal): EQN (4)\"\n);\n\n// 3. Objetivos (PT) paragraph rewrite.\nawait replaceOnce(\n  context,\n  \"Fornecer aos alunos conceitos fundamentais para compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica por meio da experimenta\u00e7\u00e3o, desenvolvendo a capacidade de realizarem pr\u00e1ticas no laborat\u00f3rio que estimulem o seu pensamento cient\u00edfico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\",\n  \"Fornecer aos alunos conceitos fundamentos para a compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica, de forma a capacit\u00e1-lo a descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\"\n);\n\n// 4. Objetivos (EN) paragraph rewrite.\nawait replaceOnce(\n  context,\n  \"Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.\",\n  \"Provide students with fundamental concepts for understanding Inorganic Chemistry, in order to enable them to describe and interpret the properties of elements and their compounds, especially those of an inorganic nature with industrial interest\"\n);\n\n// 5. Programa resumido (PT) rewrite.\nawait replaceOnce(\n  context,\n  \"Compostos de Coordena\u00e7\u00e3o. Materiais inorg\u00e2nicos de interesse industrial. Purifica\u00e7\u00e3o e Identifica\u00e7\u00e3o de Compostos Inorg\u00e2nicos. S\u00edntese de sais e obten\u00e7\u00e3o de Compostos de Alum\u00ednio.\",\n  \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Hidrog\u00eanio. Metais alcalinos. Metais alcalino terrosos. Alum\u00ednio. Metais de transi\u00e7\u00e3o. Compostos de coordena\u00e7\u00e3o. Halog\u00eanios.\"\n);\n\n// 6. Programa resumido (EN) rewrite.\nawait replaceOnce(\n  context,\n  \"Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\",\n  \"Methods for separating and obtaining elements, mineral extraction. Hydrogen. Alkaline metals. Alkaline earth metals. Aluminum. Transition metals. Coordination compounds. Halogens.\"\n);\n\n// 7. Programa (PT) rewrite.\nawait replaceOnce(\n  context,\n  \"Compostos de Coordena\u00e7\u00e3o: Estrutura, liga\u00e7\u00f5es, rea\u00e7\u00f5es e aplica\u00e7\u00f5es. Exemplos e aplica\u00e7\u00f5es de materiais inorg\u00e2nicos de interesse industrial. S\u00ednteses: Sal Simples, Sal Duplo e Sal Complexo. Prepara\u00e7\u00e3o de Compostos de Alum\u00ednio.\",\n  \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Propriedades, obten\u00e7\u00e3o e aplica\u00e7\u00f5es dos seguintes elementos/grupos e seus compostos: Hidrog\u00eanio; Metais alcalinos (ind\u00fastria cloro-\u00e1lcali; processo Solvay); Metais alcalino terrosos; Alum\u00ednio (processo Bayer); Metais de transi\u00e7\u00e3o; Compostos de coordena\u00e7\u00e3o e Halog\u00eanios.\"\n);\n\n// 8. Programa (EN) rewrite.\nawait replaceOnce(\n  context,\n  \"Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\",\n  \"Methods of separating and obtaining the elements, mineral extraction. Properties, obtaining and applications of the following elements/groups and their compounds: Hydrogen; Alkali metals (chlor-alkali industry; Solvay process); Alkaline earth metals; Aluminum (Bayer process); Transition metals; Coordination compounds and Halogens.\"\n);\n\n// 9. M\u00e9todo: drop \"e pr\u00e1ticas\" from the class format sentence.\nawait replaceOnce(\n  context,\n  \"Ser\u00e3o oferecidas aulas expositivas e pr\u00e1ticas.\",\n  \"Ser\u00e3o oferecidas aulas expositivas.\"\n);\n\n// 10. Bibliografia: prepend the WELLER reference and fix a missing space.\nawait replaceOnce(\n  context,\n  \"CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\",\n  \"WELLER, Mark; OVERTON, Tina; ROURKE, Jonathan; et al. Qu\u00edmica inorg\u00e2nica. Porto Alegre, Bookman, 6\u00aa Ed, 2017. E-book. CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman, Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $doc / $app resolve against the loaded document; ActiveDocument is\n# the live document we are editing.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdFindContinue(1) ... last arg Replace:=wdReplaceAll(2)\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1. Ativa\u00e7\u00e3o date bump.\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2022\" \"Ativa\u00e7\u00e3o: 01/01/2024\"\n\n# 2. Curso (semestre ideal): drop the EQD entry.\nReplace-Text \"Curso (semestre ideal): EQD (3), EQN (4)\" \"Curso (semestre ideal): EQN (4)\"\n\n# 3. Objetivos (PT) paragraph rewrite.\nReplace-Text \"Fornecer aos alunos conceitos fundamentais para compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica por meio da experimenta\u00e7\u00e3o, desenvolvendo a capacidade de realizarem pr\u00e1ticas no laborat\u00f3rio que estimulem o seu pensamento cient\u00edfico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\" \"Fornecer aos alunos conceitos fundamentos para a compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica, de forma a capacit\u00e1-lo a descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\"\n\n# 4. Objetivos (EN) paragraph rewrite.\nReplace-Text \"Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.\" \"Provide students with fundamental concepts for understanding Inorganic Chemistry, in order to enable them to describe and interpret the properties of elements and their compounds, especially those of an inorganic nature with industrial interest\"\n\n# 5. Programa resumido (PT) rewrite.\nReplace-Text \"Compostos de Coordena\u00e7\u00e3o. Materiais inorg\u00e2nicos de interesse industrial. Purifica\u00e7\u00e3o e Identifica\u00e7\u00e3o de Compostos Inorg\u00e2nicos. S\u00edntese de sais e obten\u00e7\u00e3o de Compostos de Alum\u00ednio.\" \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Hidrog\u00eanio. Metais alcalinos. Metais alcalino terrosos. Alum\u00ednio. Metais de transi\u00e7\u00e3o. Compostos de coordena\u00e7\u00e3o. Halog\u00eanios.\"\n\n# 6. Programa resumido (EN) rewrite.\nReplace-Text \"Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\" \"Methods for separating and obtaining elements, mineral extraction. Hydrogen. Alkaline metals. Alkaline earth metals. Aluminum. Transition metals. Coordination compounds. Halogens.\"\n\n# 7. Programa (PT) rewrite.\nReplace-Text \"Compostos de Coordena\u00e7\u00e3o: Estrutura, liga\u00e7\u00f5es, rea\u00e7\u00f5es e aplica\u00e7\u00f5es. Exemplos e aplica\u00e7\u00f5es de materiais inorg\u00e2nicos de interesse industrial. S\u00ednteses: Sal Simples, Sal Duplo e Sal Complexo. Prepara\u00e7\u00e3o de Compostos de Alum\u00ednio.\" \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Propriedades, obten\u00e7\u00e3o e aplica\u00e7\u00f5es dos seguintes elementos/grupos e seus compostos: Hidrog\u00eanio; Metais alcalinos (ind\u00fastria cloro-\u00e1lcali; processo Solvay); Metais alcalino terrosos; Alum\u00ednio (processo Bayer); Metais de transi\u00e7\u00e3o; Compostos de coordena\u00e7\u00e3o e Halog\u00eanios.\"\n\n# 8. Programa (EN) rewrite.\nReplace-Text \"Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\" \"Methods of separating and obtaining the elements, mineral extraction. Properties, obtaining and applications of the following elements/groups and their compounds: Hydrogen; Alkali metals (chlor-alkali industry; Solvay process); Alkaline earth metals; Aluminum (Bayer process); Transition metals; Coordination compounds and Halogens.\"\n\n# 9. M\u00e9todo: drop \"e pr\u00e1ticas\" from the class format sentence.\nReplace-Text \"Ser\u00e3o oferecidas aulas expositivas e pr\u00e1ticas.\" \"Ser\u00e3o oferecidas aulas expositivas.\"\n\n# 10. Bibliografia: prepend the WELLER reference and fix a missing space.\nReplace-Text \"CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\" \"WELLER, Mark; OVERTON, Tina; ROURKE, Jonathan; et al. Qu\u00edmica inorg\u00e2nica. Porto Alegre, Bookman, 6\u00aa Ed, 2017. E-book. CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman, Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\"\n"}
